$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44322
$ws.Range("M2").Value = 60
$ws.Range("Q2").Value = "$/caja 10 kilos empedrada"
$ws.Range("S2").Value = 11500
$ws.Range("T2").Value = 1

# Row 4
$ws.Range("D4").Value = 44327
$ws.Range("M4").Value = 60

# Row 5
$ws.Range("D5").Value = 44323
$ws.Range("M5").Value = 80

# Row 6
$ws.Range("D6").Value = 44316

# Row 8
$ws.Range("D8").Value = 44309
$ws.Range("Q8").Value = "$/caja 14 kilos granel"
$ws.Range("S8").Value = 821
$ws.Range("T8").Value = 14

# Row 9
$ws.Range("D9").Value = 44306
$ws.Range("M9").Value = 80

# Row 10
$ws.Range("D10").Value = 44313
$ws.Range("M10").Value = 120
